# Fruta / hortaliza, semanal
#
# Weekly refresh of the "Arandano (blue)" price sheet for Terminal
# Hortofruticola Agro Chillan: two new daily observations for 2023-02-22
# (serial 44979) are inserted at the top of the data block (row 9),
# pushing every existing record down by two rows (old row 9 -> row 11,
# ..., old row 34 -> row 36).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right above the current first data row so the
# existing records (9..34) shift down to (11..36).
$ws.Rows("9:10").Insert()

# --- New row 9: Primera, 2023-02-22 ------------------------------------
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C9").Value = "Ñuble"
$ws.Range("D9").Value = 44979
$ws.Range("E9").Value = 16
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100101
$ws.Range("H9").Value = "Berries"
$ws.Range("I9").Value = 100101001
$ws.Range("J9").Value = "Arándano (blue)"
$ws.Range("K9").Value = "Sin especificar"
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 30
$ws.Range("N9").Value = 3000
$ws.Range("O9").Value = 3000
$ws.Range("P9").Value = 3000
$ws.Range("Q9").Value = "$/bandeja 2 kilos"
$ws.Range("R9").Value = "Provincia de Diguillín"
$ws.Range("S9").Value = 1500
$ws.Range("T9").Value = 2

# --- New row 10: Segunda, 2023-02-22 -----------------------------------
$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C10").Value = "Ñuble"
$ws.Range("D10").Value = 44979
$ws.Range("E10").Value = 16
$ws.Range("F10").Value = "Fruta"
$ws.Range("G10").Value = 100101
$ws.Range("H10").Value = "Berries"
$ws.Range("I10").Value = 100101001
$ws.Range("J10").Value = "Arándano (blue)"
$ws.Range("K10").Value = "Sin especificar"
$ws.Range("L10").Value = "Segunda"
$ws.Range("M10").Value = 30
$ws.Range("N10").Value = 2500
$ws.Range("O10").Value = 2500
$ws.Range("P10").Value = 2500
$ws.Range("Q10").Value = "$/bandeja 2 kilos"
$ws.Range("R10").Value = "Provincia de Diguillín"
$ws.Range("S10").Value = 1250
$ws.Range("T10").Value = 2

Write-Output "done"
